# Daily attendance processing - 2025-11-06 17:22:30
# Normalizes the "Recorded By" column (G) so that entries whose list of
# recorders starts with "System"/"system" have that token swapped with the
# last entry in the comma-separated list instead of leading it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ", "

    if ($parts.Count -ge 2 -and $parts[0].Trim().ToLower() -eq "system") {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
        $cell.Value = [string]::Join(", ", $parts)
    }
}
